# Updates odds & correct-score prices for Jogos_da_Semana_FlashScore_2025-02-10.xlsx
# (values refreshed per the latest FlashScore odds snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("G4").Value = 2.35
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 3.3
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 4.33
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 6
$ws.Range("T4").Value = 1.13
$ws.Range("U4").Value = 1.58
$ws.Range("V4").Value = 2.2
$ws.Range("W4").Value = 2.25
$ws.Range("X4").Value = 1.54
$ws.Range("Z4").Value = 9.5
$ws.Range("AE4").Value = 6
$ws.Range("AG4").Value = 21
$ws.Range("AH4").Value = 81
$ws.Range("AM4").Value = 41
$ws.Range("AO4").Value = 51
$ws.Range("AP4").Value = 2.05
$ws.Range("AQ4").Value = 1.8

# Row 5
$ws.Range("U5").Value = 1.54
$ws.Range("X5").Value = 1.54
$ws.Range("AP5").Value = 1.92
$ws.Range("AQ5").Value = 1.82

# Row 7
$ws.Range("G7").Value = 1.7
$ws.Range("I7").Value = 4.9
$ws.Range("J7").Value = 2.35
$ws.Range("L7").Value = 5.3
$ws.Range("Q7").Value = 2.12
$ws.Range("R7").Value = 1.57
$ws.Range("Z7").Value = 6.9
$ws.Range("AB7").Value = 12.5
$ws.Range("AC7").Value = 15.5
$ws.Range("AF7").Value = 6.7
$ws.Range("AM7").Value = 100

# Row 8
$ws.Range("G8").Value = 1.67
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 5.75
$ws.Range("J8").Value = 2.4
$ws.Range("K8").Value = 1.95
$ws.Range("L8").Value = 6.5
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 2.5
$ws.Range("U8").Value = 1.54
$ws.Range("X8").Value = 1.47
$ws.Range("Y8").Value = 4.75
$ws.Range("Z8").Value = 6.5
$ws.Range("AB8").Value = 12
$ws.Range("AG8").Value = 26
$ws.Range("AJ8").Value = 10
$ws.Range("AK8").Value = 26
$ws.Range("AL8").Value = 21
$ws.Range("AM8").Value = 67
$ws.Range("AO8").Value = 67
$ws.Range("AP8").Value = 1.94
$ws.Range("AQ8").Value = 1.79

# Row 9
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.07
$ws.Range("R9").Value = 1.69

# Row 10
$ws.Range("R10").Value = 1.44
$ws.Range("AP10").Value = 1.94
$ws.Range("AQ10").Value = 1.79
$ws.Range("AR10").Value = 4.3
$ws.Range("AS10").Value = 1.22

# Row 11
$ws.Range("R11").Value = 1.41
$ws.Range("AP11").Value = 1.99
$ws.Range("AQ11").Value = 1.74

# Row 12
$ws.Range("W12").Value = 1.47

# Row 15
$ws.Range("W15").Value = 1.87
$ws.Range("X15").Value = 1.77

# Row 16
$ws.Range("G16").Value = 1.67
$ws.Range("H16").Value = 3.8
$ws.Range("J16").Value = 2.25
$ws.Range("M16").Value = 1.04
$ws.Range("N16").Value = 13
$ws.Range("Q16").Value = 1.83
$ws.Range("R16").Value = 2.03
$ws.Range("W16").Value = 1.77
$ws.Range("X16").Value = 1.92
$ws.Range("Z16").Value = 8
$ws.Range("AI16").Value = 201
$ws.Range("AJ16").Value = 15
$ws.Range("AL16").Value = 17

# Row 17
$ws.Range("N17").Value = 21
$ws.Range("S17").Value = 2
$ws.Range("T17").Value = 1.73

# Row 18
$ws.Range("G18").Value = 4.75
$ws.Range("H18").Value = 4.5
$ws.Range("I18").Value = 1.6
$ws.Range("K18").Value = 2.62
$ws.Range("L18").Value = 2.1
$ws.Range("O18").Value = 1.11
$ws.Range("P18").Value = 6.5
$ws.Range("Q18").Value = 1.4
$ws.Range("R18").Value = 2.88
$ws.Range("S18").Value = 1.91
$ws.Range("T18").Value = 1.8
$ws.Range("U18").Value = 1.22
$ws.Range("V18").Value = 4
$ws.Range("W18").Value = 1.44
$ws.Range("X18").Value = 2.63
$ws.Range("AC18").Value = 29
$ws.Range("AE18").Value = 23
$ws.Range("AF18").Value = 9.5
$ws.Range("AG18").Value = 13
$ws.Range("AI18").Value = 81
$ws.Range("AN18").Value = 11
$ws.Range("AO18").Value = 17

# Row 19
$ws.Range("K19").Value = 2.37
$ws.Range("N19").Value = 17
$ws.Range("O19").Value = 1.14
$ws.Range("P19").Value = 5.5
$ws.Range("Q19").Value = 1.53
$ws.Range("R19").Value = 2.4
$ws.Range("S19").Value = 2.25
$ws.Range("T19").Value = 1.57
$ws.Range("AR19").Value = 1.79
$ws.Range("AS19").Value = 1.94

# Row 21
$ws.Range("Q21").Value = 1.92
$ws.Range("R21").Value = 1.82

# Row 23
$ws.Range("M23").Value = 1.07
$ws.Range("N23").Value = 9
$ws.Range("Q23").Value = 2.2
$ws.Range("R23").Value = 1.65

# Row 24
$ws.Range("U24").Value = 1.47

# Row 25
$ws.Range("G25").Value = 1.85
$ws.Range("H25").Value = 3.75
$ws.Range("I25").Value = 3.7
$ws.Range("J25").Value = 2.37
$ws.Range("K25").Value = 2.27
$ws.Range("L25").Value = 4
$ws.Range("M25").Value = 1.04
$ws.Range("N25").Value = 8.75
$ws.Range("O25").Value = 1.2
$ws.Range("P25").Value = 4.05
$ws.Range("Q25").Value = 1.6
$ws.Range("R25").Value = 2.2
$ws.Range("S25").Value = 2.42
$ws.Range("T25").Value = 1.5
$ws.Range("U25").Value = 1.32
$ws.Range("V25").Value = 3.1
$ws.Range("W25").Value = 1.57
$ws.Range("X25").Value = 2.25
$ws.Range("Y25").Value = 9.5
$ws.Range("Z25").Value = 10.25
$ws.Range("AA25").Value = 8.25
$ws.Range("AB25").Value = 16.5
$ws.Range("AC25").Value = 13
$ws.Range("AD25").Value = 20
$ws.Range("AE25").Value = 8.75
$ws.Range("AF25").Value = 7.5
$ws.Range("AG25").Value = 12.5
$ws.Range("AH25").Value = 45
$ws.Range("AI25").Value = 250
$ws.Range("AJ25").Value = 14
$ws.Range("AK25").Value = 23
$ws.Range("AL25").Value = 12.5
$ws.Range("AM25").Value = 55
$ws.Range("AN25").Value = 29
$ws.Range("AO25").Value = 30

# Row 26
$ws.Range("G26").Value = 2.9
$ws.Range("I26").Value = 2.3
$ws.Range("U26").Value = 1.41
$ws.Range("V26").Value = 2.62
$ws.Range("Y26").Value = 9.5
$ws.Range("AK26").Value = 11
$ws.Range("AN26").Value = 19

# Row 27
$ws.Range("U27").Value = 1.27
